# Apply "hybrid bold + color" highlighting to quantitative metrics inside
# specific bullet paragraphs, matching the supplied OOXML diff. Each target
# paragraph currently holds one run with the full sentence text; we rebuild
# it as alternating plain / bold+colored runs without disturbing any other
# paragraph in the document.

$d = $word.ActiveDocument

# BGR-packed integer for hex color 2C3E50 (Word's Font.Color / wdColor is
# stored internally as 0x00BBGGRR).
$HighlightColorBGR = 0x50 * 65536 + 0x3E * 256 + 0x2C

# Each entry: the exact current paragraph text (sans trailing paragraph
# mark), and the ordered list of (Bold, Text) segments that should replace
# it - concatenating the segment Text values reproduces the original text
# exactly.
$targets = @(
    @{
        Match = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
        Segments = @(
            @{ Bold = $false; Text = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from " },
            @{ Bold = $true;  Text = "23%" },
            @{ Bold = $false; Text = " to " },
            @{ Bold = $true;  Text = "64%" }
        )
    },
    @{
        Match = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
        Segments = @(
            @{ Bold = $false; Text = "• Achieved " },
            @{ Bold = $true;  Text = "87%" },
            @{ Bold = $false; Text = " prediction accuracy for voter turnout vs. industry standard of " },
            @{ Bold = $true;  Text = "71%" },
            @{ Bold = $false; Text = ", reducing polling error margins from " },
            @{ Bold = $true;  Text = "±4.2%" },
            @{ Bold = $false; Text = " to " },
            @{ Bold = $true;  Text = "±2.1%" }
        )
    },
    @{
        Match = "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
        Segments = @(
            @{ Bold = $false; Text = "• Wrote RFP and analyzed bids from " },
            @{ Bold = $true;  Text = "1,200" },
            @{ Bold = $false; Text = " vendors for research platform development" }
        )
    },
    @{
        Match = "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"
        Segments = @(
            @{ Bold = $false; Text = "• Created comprehensive meta-analysis framework handling millions of survey responses that became the " },
            @{ Bold = $true;  Text = "`$400M" },
            @{ Bold = $false; Text = " Polling Consortium Database at The Analyst Institute, now valued at " },
            @{ Bold = $true;  Text = "`$1B" },
            @{ Bold = $false; Text = "+" }
        )
    },
    @{
        Match = "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M"
        Segments = @(
            @{ Bold = $false; Text = "• Algorithm reduced mapping costs by " },
            @{ Bold = $true;  Text = "73.5%" },
            @{ Bold = $false; Text = ", saving campaigns and organizations " },
            @{ Bold = $true;  Text = "`$4.7M" }
        )
    },
    @{
        Match = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
        Segments = @(
            @{ Bold = $false; Text = "• Achieved " },
            @{ Bold = $true;  Text = "87%" },
            @{ Bold = $false; Text = " prediction accuracy for voter turnout vs. industry standard of " },
            @{ Bold = $true;  Text = "71%" }
        )
    }
)

foreach ($target in $targets) {
    $matchText = $target.Match
    $paraFound = $null
    foreach ($p in $d.Paragraphs) {
        $pText = $p.Range.Text
        # Paragraph Range.Text always carries the trailing paragraph mark
        # (CR, char 13); trim it before comparing against the diff text.
        if ($pText.Length -gt 0 -and $pText.Substring($pText.Length - 1) -eq "`r") {
            $pText = $pText.Substring(0, $pText.Length - 1)
        }
        if ($pText -eq $matchText) {
            $paraFound = $p
            break
        }
    }
    if ($paraFound -eq $null) {
        Write-Output "WARNING: paragraph not found for: $matchText"
        continue
    }

    $pRange = $paraFound.Range
    $pStart = $pRange.Start

    # Replace the paragraph's whole text (minus paragraph mark) with the
    # first segment, then insert each subsequent segment right after the
    # growing content, applying bold/color per-segment.
    $first = $target.Segments[0]
    $pRange.Text = $first.Text
    $firstRange = $d.Range($pStart, $pStart + $first.Text.Length)
    $firstRange.Font.Bold = $first.Bold
    if ($first.Bold) {
        $firstRange.Font.Color = $HighlightColorBGR
    }

    $pos = $pStart + $first.Text.Length

    for ($i = 1; $i -lt $target.Segments.Count; $i++) {
        $seg = $target.Segments[$i]
        $insertPoint = $d.Range($pos, $pos)
        $insertPoint.InsertAfter($seg.Text)
        $segRange = $d.Range($pos, $pos + $seg.Text.Length)
        $segRange.Font.Bold = $seg.Bold
        if ($seg.Bold) {
            $segRange.Font.Color = $HighlightColorBGR
        }
        $pos = $pos + $seg.Text.Length
    }

    Write-Output "Updated paragraph: $matchText"
}

Write-Output "Done"
